$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H1) mirrors the existing header styling (bold,
# bordered, centered) by copying the format from the neighboring "sum"
# header cell (G1), then overwriting just the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column on row 2.
$ws.Range("H2").Value = 1
